$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 73303.57209245772
$ws.Cells.Item(3, 2).Value = 69996.49426258687
$ws.Cells.Item(4, 2).Value = 67352.29750351419
$ws.Cells.Item(5, 2).Value = 64777.80292230453
$ws.Cells.Item(6, 2).Value = 64645.22173689886
$ws.Cells.Item(7, 2).Value = 67323.31236904059
$ws.Cells.Item(8, 2).Value = 67299.50380321605
$ws.Cells.Item(9, 2).Value = 72662.17442876496
$ws.Cells.Item(10, 2).Value = 88486.30291684365
$ws.Cells.Item(11, 2).Value = 97673.19597505673
$ws.Cells.Item(12, 2).Value = 102152.0786721448
$ws.Cells.Item(13, 2).Value = 102207.5360816811
$ws.Cells.Item(14, 2).Value = 102787.0150887881
$ws.Cells.Item(15, 2).Value = 106750.4892649202
$ws.Cells.Item(16, 2).Value = 106887.273802493
$ws.Cells.Item(17, 2).Value = 104380.0409036494
$ws.Cells.Item(18, 2).Value = 100139.0962381116
$ws.Cells.Item(19, 2).Value = 91146.88819159759
$ws.Cells.Item(20, 2).Value = 90356.49322132593
$ws.Cells.Item(21, 2).Value = 88933.55142855708
$ws.Cells.Item(22, 2).Value = 87187.11569680319
$ws.Cells.Item(23, 2).Value = 84791.4262532634
$ws.Cells.Item(24, 2).Value = 81293.23952288678
$ws.Cells.Item(25, 2).Value = 76753.33483251851
